$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Data" (G) and "nrworkshoop" (H) values for several rows.
# (Row 14 previously held the duplicate date string " 11/8/2022", which
# becomes unused and is dropped from the shared strings table.)

$ws.Range("G6").Value  = "  8/8/2022"
$ws.Range("H6").Value  = 1

$ws.Range("G7").Value  = " 8/11/2022"
$ws.Range("H7").Value  = 2

$ws.Range("G8").Value  = "  9/8/2022"
$ws.Range("H8").Value  = 1

$ws.Range("G10").Value = "  9/8/2022"
$ws.Range("H10").Value = 1

$ws.Range("G12").Value = " 8/11/2022"
$ws.Range("H12").Value = 2

$ws.Range("G13").Value = " 11/9/2022"
$ws.Range("H13").Value = 2

$ws.Range("G14").Value = " 11/9/2022"

$ws.Range("G15").Value = "  8/8/2022"
$ws.Range("H15").Value = 1

$ws.Range("G16").Value = " 11/9/2022"
$ws.Range("H16").Value = 2

$ws.Range("G19").Value = "  9/8/2022"
$ws.Range("H19").Value = 1

$ws.Range("G21").Value = " 11/9/2022"
$ws.Range("H21").Value = 2

$ws.Range("G23").Value = "  9/8/2022"
$ws.Range("H23").Value = 1

$ws.Range("G25").Value = "  9/8/2022"
$ws.Range("H25").Value = 1

$ws.Range("G26").Value = "  9/8/2022"
$ws.Range("H26").Value = 1
